$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-type formatting on columns that hold numeric-looking or date-looking
# text (Antal, Startdatum, Slutdatum) so Excel does not auto-convert the literal
# strings we are about to write into numbers/dates.
$ws.Range("I3:I9").NumberFormat = "@"
$ws.Range("Y3:Y9").NumberFormat = "@"
$ws.Range("AA3:AA9").NumberFormat = "@"

# New row 3 gets the content that was in old row 4
$ws.Range("A3").Value = 3056139
$ws.Range("B3").Value = 96334
$ws.Range("C3").Value = "Godkänd baserat på observatörens uppgifter"
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 220787
$ws.Range("F3").Value = "Knärot"
$ws.Range("G3").Value = "Goodyera repens"
$ws.Range("H3").Value = "(L.) R. Br."
$ws.Range("I3").Value = "215"
$ws.Range("J3").Value = "plantor/tuvor"
$ws.Range("P3").Value = "Östtomten, 360 m NÖ om, Vrm"
$ws.Range("Q3").Value = 365559.0432148376
$ws.Range("R3").Value = 6618565.766507677
$ws.Range("S3").Value = 10
$ws.Range("T3").Value = "Värmland"
$ws.Range("U3").Value = "Arvika"
$ws.Range("V3").Value = "Värmland"
$ws.Range("W3").Value = "Arvika"
$ws.Range("Y3").Value = "2011-10-23"
$ws.Range("Z3").Value = "00:00"
$ws.Range("AA3").Value = "2011-10-23"
$ws.Range("AB3").Value = "00:00"
$ws.Range("AC3").Value = "215 bladrosetter med 32 blomstänglar"
$ws.Range("AD3").Value = $False
$ws.Range("AE3").Value = $False
$ws.Range("AG3").Value = $False
$ws.Range("AH3").Value = "Granskog"
$ws.Range("AI3").Value = "äldre granskog med inslag av tall"
$ws.Range("AW3").Value = "Per Larsson"
$ws.Range("AX3").Value = "Per Larsson, Elvi Eriksson"

# New row 4 gets the content that was in old row 5
$ws.Range("A4").Value = 3056366
$ws.Range("B4").Value = 96334
$ws.Range("C4").Value = "Godkänd baserat på observatörens uppgifter"
$ws.Range("D4").Value = "VU"
$ws.Range("E4").Value = 220787
$ws.Range("F4").Value = "Knärot"
$ws.Range("G4").Value = "Goodyera repens"
$ws.Range("H4").Value = "(L.) R. Br."
$ws.Range("I4").Value = "600"
$ws.Range("J4").Value = "plantor/tuvor"
$ws.Range("P4").Value = "Östtomten, 320 m NNÖ om, Vrm"
$ws.Range("Q4").Value = 365463.5963024567
$ws.Range("R4").Value = 6618600.539037504
$ws.Range("S4").Value = 10
$ws.Range("T4").Value = "Värmland"
$ws.Range("U4").Value = "Arvika"
$ws.Range("V4").Value = "Värmland"
$ws.Range("W4").Value = "Arvika"
$ws.Range("Y4").Value = "2011-10-23"
$ws.Range("Z4").Value = "00:00"
$ws.Range("AA4").Value = "2011-10-23"
$ws.Range("AB4").Value = "00:00"
$ws.Range("AC4").Value = "mittkoordinat, 490 bladrosetter-41 blomstänglar 525-075, 110 bladrosetter-7 blomstänglar 534-088"
$ws.Range("AD4").Value = $False
$ws.Range("AE4").Value = $False
$ws.Range("AG4").Value = $False
$ws.Range("AH4").Value = "Granskog"
$ws.Range("AI4").Value = "äldre granskog med inslag av tall"
$ws.Range("AW4").Value = "Per Larsson"
$ws.Range("AX4").Value = "Per Larsson, Elvi Eriksson"

# New row 5 gets the content that was in old row 6
$ws.Range("A5").Value = 3056367
$ws.Range("B5").Value = 96334
$ws.Range("C5").Value = "Godkänd baserat på observatörens uppgifter"
$ws.Range("D5").Value = "VU"
$ws.Range("E5").Value = 220787
$ws.Range("F5").Value = "Knärot"
$ws.Range("G5").Value = "Goodyera repens"
$ws.Range("H5").Value = "(L.) R. Br."
$ws.Range("I5").Value = "160"
$ws.Range("J5").Value = "plantor/tuvor"
$ws.Range("P5").Value = "Östtomten, 270 m NNÖ om, Vrm"
$ws.Range("Q5").Value = 365428.9013190148
$ws.Range("R5").Value = 6618579.045586957
$ws.Range("S5").Value = 10
$ws.Range("T5").Value = "Värmland"
$ws.Range("U5").Value = "Arvika"
$ws.Range("V5").Value = "Värmland"
$ws.Range("W5").Value = "Arvika"
$ws.Range("Y5").Value = "2011-10-23"
$ws.Range("Z5").Value = "00:00"
$ws.Range("AA5").Value = "2011-10-23"
$ws.Range("AB5").Value = "00:00"
$ws.Range("AC5").Value = "160 bladrosetter med 14 blomstänglar"
$ws.Range("AD5").Value = $False
$ws.Range("AE5").Value = $False
$ws.Range("AG5").Value = $False
$ws.Range("AH5").Value = "Granskog"
$ws.Range("AI5").Value = "äldre granskog med inslag av tall"
$ws.Range("AW5").Value = "Per Larsson"
$ws.Range("AX5").Value = "Per Larsson, Elvi Eriksson"

# New row 6 gets the content that was in old row 7
$ws.Range("A6").Value = 3056136
$ws.Range("B6").Value = 96334
$ws.Range("C6").Value = "Godkänd baserat på observatörens uppgifter"
$ws.Range("D6").Value = "VU"
$ws.Range("E6").Value = 220787
$ws.Range("F6").Value = "Knärot"
$ws.Range("G6").Value = "Goodyera repens"
$ws.Range("H6").Value = "(L.) R. Br."
$ws.Range("I6").Value = "80"
$ws.Range("J6").Value = "plantor/tuvor"
$ws.Range("P6").Value = "Östtomten, 490 m NÖ om, Vrm"
$ws.Range("Q6").Value = 365699.2698654335
$ws.Range("R6").Value = 6618593.572883067
$ws.Range("S6").Value = 10
$ws.Range("T6").Value = "Värmland"
$ws.Range("U6").Value = "Arvika"
$ws.Range("V6").Value = "Värmland"
$ws.Range("W6").Value = "Arvika"
$ws.Range("Y6").Value = "2011-10-23"
$ws.Range("Z6").Value = "00:00"
$ws.Range("AA6").Value = "2011-10-23"
$ws.Range("AB6").Value = "00:00"
$ws.Range("AC6").Value = "80 bladrosetter med 12 blomstänglar"
$ws.Range("AD6").Value = $False
$ws.Range("AE6").Value = $False
$ws.Range("AG6").Value = $False
$ws.Range("AH6").Value = "Granskog"
$ws.Range("AI6").Value = "äldre granskog med inslag av tall"
$ws.Range("AW6").Value = "Per Larsson"
$ws.Range("AX6").Value = "Per Larsson, Elvi Eriksson"

# New row 7 gets the content that was in old row 9
$ws.Range("A7").Value = 3056138
$ws.Range("B7").Value = 96334
$ws.Range("C7").Value = "Godkänd baserat på observatörens uppgifter"
$ws.Range("D7").Value = "VU"
$ws.Range("E7").Value = 220787
$ws.Range("F7").Value = "Knärot"
$ws.Range("G7").Value = "Goodyera repens"
$ws.Range("H7").Value = "(L.) R. Br."
$ws.Range("I7").Value = "36"
$ws.Range("J7").Value = "plantor/tuvor"
$ws.Range("P7").Value = "Östtomten, 330 m NÖ om, Vrm"
$ws.Range("Q7").Value = 365557.7136587457
$ws.Range("R7").Value = 6618514.767422446
$ws.Range("S7").Value = 25
$ws.Range("T7").Value = "Värmland"
$ws.Range("U7").Value = "Arvika"
$ws.Range("V7").Value = "Värmland"
$ws.Range("W7").Value = "Arvika"
$ws.Range("Y7").Value = "2011-10-23"
$ws.Range("Z7").Value = "00:00"
$ws.Range("AA7").Value = "2011-10-23"
$ws.Range("AB7").Value = "00:00"
$ws.Range("AC7").Value = "mittkoordinat, 1 bladrosett 426-167, 20 bladrosetter med 3 blomstänglar 448-177, 15 bladrosetter med 1 blomstängel 461-180"
$ws.Range("AD7").Value = $False
$ws.Range("AE7").Value = $False
$ws.Range("AG7").Value = $False
$ws.Range("AH7").Value = "Granskog"
$ws.Range("AI7").Value = "äldre granskog med inslag av tall"
$ws.Range("AW7").Value = "Per Larsson"
$ws.Range("AX7").Value = "Per Larsson, Elvi Eriksson"

# New row 8 gets the content that was in old row 3
$ws.Range("A8").Value = 3057790
$ws.Range("B8").Value = 96334
$ws.Range("C8").Value = "Godkänd baserat på observatörens uppgifter"
$ws.Range("D8").Value = "VU"
$ws.Range("E8").Value = 220787
$ws.Range("F8").Value = "Knärot"
$ws.Range("G8").Value = "Goodyera repens"
$ws.Range("H8").Value = "(L.) R. Br."
$ws.Range("I8").Value = "30"
$ws.Range("J8").Value = "plantor/tuvor"
$ws.Range("P8").Value = "Gärdet, 500 m ÖNÖ om, Vrm"
$ws.Range("Q8").Value = 365466.3266593472
$ws.Range("R8").Value = 6618760.65168799
$ws.Range("S8").Value = 25
$ws.Range("T8").Value = "Värmland"
$ws.Range("U8").Value = "Arvika"
$ws.Range("V8").Value = "Värmland"
$ws.Range("W8").Value = "Arvika"
$ws.Range("Y8").Value = "2011-12-08"
$ws.Range("Z8").Value = "00:00"
$ws.Range("AA8").Value = "2011-12-08"
$ws.Range("AB8").Value = "00:00"
$ws.Range("AC8").Value = "30 bladrosetter, 2 smågrupper"
$ws.Range("AD8").Value = $False
$ws.Range("AE8").Value = $False
$ws.Range("AG8").Value = $False
$ws.Range("AH8").Value = "Granskog"
$ws.Range("AI8").Value = "granskog med en del tallar, gallrad för ett tag sedan"
$ws.Range("AW8").Value = "Per Larsson"
$ws.Range("AX8").Value = "Per Larsson"

# New row 9 gets the content that was in old row 8
$ws.Range("A9").Value = 3056137
$ws.Range("B9").Value = 96334
$ws.Range("C9").Value = "Godkänd baserat på observatörens uppgifter"
$ws.Range("D9").Value = "VU"
$ws.Range("E9").Value = 220787
$ws.Range("F9").Value = "Knärot"
$ws.Range("G9").Value = "Goodyera repens"
$ws.Range("H9").Value = "(L.) R. Br."
$ws.Range("I9").Value = "110"
$ws.Range("J9").Value = "plantor/tuvor"
$ws.Range("P9").Value = "Östtomten, 350 m NÖ om, Vrm"
$ws.Range("Q9").Value = 365594.945751337
$ws.Range("R9").Value = 6618494.221399919
$ws.Range("S9").Value = 25
$ws.Range("T9").Value = "Värmland"
$ws.Range("U9").Value = "Arvika"
$ws.Range("V9").Value = "Värmland"
$ws.Range("W9").Value = "Arvika"
$ws.Range("Y9").Value = "2011-10-23"
$ws.Range("Z9").Value = "00:00"
$ws.Range("AA9").Value = "2011-10-23"
$ws.Range("AB9").Value = "00:00"
$ws.Range("AC9").Value = "mittkoordinat, 100 bladrosetter med 3 blomstänglar 410-210, 10 bladrosetter med 1 blomstängel 435-209"
$ws.Range("AD9").Value = $False
$ws.Range("AE9").Value = $False
$ws.Range("AG9").Value = $False
$ws.Range("AH9").Value = "Granskog"
$ws.Range("AI9").Value = "äldre granskog med inslag av tall"
$ws.Range("AW9").Value = "Per Larsson"
$ws.Range("AX9").Value = "Per Larsson, Elvi Eriksson"

# Restore the default ("Normal") cell style on the text-forced columns so the
# worksheet XML does not carry a lingering explicit style index, matching the
# original workbook formatting (values remain stored as text).
$ws.Range("I3:I9").Style = "Normal"
$ws.Range("Y3:Y9").Style = "Normal"
$ws.Range("AA3:AA9").Style = "Normal"

"rows 3-9 rearranged"